$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Order mirrors the sheet's row order (2..51).
# Values are stored as inline-string text in the workbook (prices/links/
# percent deltas are text, not numbers), so each cell is forced to Text
# number-format before assignment and then returned to the default
# "Normal" style so no stray formatting is left behind.
$updates = [ordered]@{
    'D2'  = '22.451.13'
    'E2'  = '  +0.64%  '
    'D3'  = '1.569.70'
    'E3'  = '  +0.27%  '
    'D4'  = '1.002'
    'E4'  = '  -0.40%  '
    'E5'  = '  -0.16%  '
    'D6'  = '290.15'
    'E6'  = '  +0.13%  '
    'D7'  = '0.3689'
    'E7'  = '  -1.36%  '
    'D8'  = '49.93'
    'E8'  = '  +1.70%  '
    'D9'  = '0.3380'
    'E9'  = '  -0.13%  '
    'D10' = '1.147'
    'E10' = '  +2.00%  '
    'D11' = '0.07537'
    'E11' = '  +0.46%  '
    'E12' = '  -0.42%  '
    'D13' = '21.17'
    'E13' = '  +1.99%  '
    'D14' = '6.033'
    'E14' = '  +2.49%  '
    'D15' = '6.976'
    'E15' = '  +1.76%  '
    'D16' = '1.569.29'
    'E16' = '  +0.18%  '
    'D17' = '0.00001122'
    'E17' = '  +0.74%  '
    'D18' = '90.38'
    'E18' = '  +0.87%  '
    'D19' = '0.06782'
    'E19' = '  +0.72%  '
    'E20' = '  -0.27%  '
    'D21' = '6.364'
    'E21' = '  +3.25%  '
    'D22' = '16.40'
    'E22' = '  +0.48%  '
    'E23' = '  +3.27%  '
    'D24' = '22.457.04'
    'E24' = '  +0.67%  '
    'D25' = '2.378'
    'E25' = '  -0.11%  '
    'D26' = '2.653'
    'E26' = '  -1.81%  '
    'D27' = '19.99'
    'E27' = '  +0.08%  '
    'D28' = '149.03'
    'E28' = '  +1.10%  '
    'D29' = '5.060'
    'E29' = '  +1.33%  '
    'D30' = '124.91'
    'E30' = '  -0.08%  '
    'D31' = '1.749.17'
    'E31' = '  +0.59%  '
    'D32' = '1.064'
    'E32' = '  +8.59%  '
    'D33' = '6.209'
    'E33' = '  +4.83%  '
    'D34' = '2.014'
    'E34' = '  -0.02%  '
    'D35' = '9.805'
    'E35' = '  +0.01%  '
    'D36' = '0.08349'
    'E36' = '  -1.01%  '
    'D37' = '0.02476'
    'E37' = '  +1.59%  '
    'D38' = '1.356'
    'E38' = '  -3.22%  '
    'D39' = '0.2298'
    'E39' = '  +1.46%  '
    'D40' = '0.06544'
    'E40' = '  +2.52%  '
    'D41' = '5.405'
    'E41' = '  +1.23%  '
    'D42' = '11.25'
    'E42' = '  +3.36%  '
    'D43' = '0.6235'
    'E43' = '  +0.31%  '
    # Rows 44/45 swap places: EnergySwap <-> Frax (name, link, price, delta).
    'B44' = 'Frax'
    'C44' = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
    'D44' = '1.002'
    'E44' = '  -0.18%  '
    'B45' = 'EnergySwap'
    'C45' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D45' = '14.05'
    'E45' = '  +1.44%  '
    'D46' = '3.798'
    'E46' = '  +0.27%  '
    'D47' = '0.5867'
    'E47' = '  +1.26%  '
    'D48' = '2.071'
    'E48' = '  +1.50%  '
    'D49' = '127.73'
    'E49' = '  +3.01%  '
    'D50' = '1.249'
    'E50' = '  +0.29%  '
    'D51' = '0.07301'
    'E51' = '  -0.16%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so values like "1.002" or "290.15" don't get
    # auto-coerced into numbers by Excel's input parser.
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$addr]
    # Drop back to the default style so no extra per-cell formatting
    # (beyond the value itself) is introduced versus the original file.
    $cell.Style = 'Normal'
}
